# Update NATMI LR-pair (Vtn-Itga5) stats with refreshed TPM-derived values.
# Ligand stats (G/H/I/J) are keyed by Sending cluster (rows 2-4 = ECs, 5-7 = FAPs, 8-10 = MuSCs);
# Receptor stats (M/N/O/P) are keyed by Target cluster (D column); Edge stats (Q/R/S/T) are
# derived from the Ligand/Receptor values. All ten data rows are updated below to match the
# new pipeline run's output.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 12.056684
$ws.Range("H2").Value = 36.170052
$ws.Range("I2").Value = 0.06307822458376462
$ws.Range("J2").Value = 0.06307822458376462
$ws.Range("M2").Value = 11.146846
$ws.Range("N2").Value = 33.440538
$ws.Range("O2").Value = 0.2594806085672136
$ws.Range("P2").Value = 0.2594806085672136
$ws.Range("Q2").Value = 134.393999818664
$ws.Range("R2").Value = 1209.545998367976
$ws.Range("S2").Value = 0.01636757610233461
$ws.Range("T2").Value = 0.01636757610233461

$ws.Range("G3").Value = 12.056684
$ws.Range("H3").Value = 36.170052
$ws.Range("I3").Value = 0.06307822458376462
$ws.Range("J3").Value = 0.06307822458376462
$ws.Range("O3").Value = 0.6444737471070977
$ws.Range("P3").Value = 0.6444737471070977
$ws.Range("Q3").Value = 333.7952887119479
$ws.Range("R3").Value = 3004.157598407532
$ws.Range("S3").Value = 0.04065225975836183
$ws.Range("T3").Value = 0.04065225975836183

$ws.Range("G4").Value = 12.056684
$ws.Range("H4").Value = 36.170052
$ws.Range("I4").Value = 0.06307822458376462
$ws.Range("J4").Value = 0.06307822458376462
$ws.Range("O4").Value = 0.09604564432568881
$ws.Range("P4").Value = 0.09604564432568881
$ws.Range("Q4").Value = 49.74536778437734
$ws.Range("R4").Value = 447.708310059396
$ws.Range("S4").Value = 0.006058388723068177
$ws.Range("T4").Value = 0.006058388723068177

$ws.Range("I5").Value = 0.1315309049843414
$ws.Range("J5").Value = 0.1315309049843414
$ws.Range("M5").Value = 11.146846
$ws.Range("N5").Value = 33.440538
$ws.Range("O5").Value = 0.2594806085672136
$ws.Range("P5").Value = 0.2594806085672136
$ws.Range("Q5").Value = 280.2387755403641
$ws.Range("R5").Value = 2522.148979863277
$ws.Range("S5").Value = 0.03412971927073324
$ws.Range("T5").Value = 0.03412971927073324

$ws.Range("I6").Value = 0.1315309049843414
$ws.Range("J6").Value = 0.1315309049843414
$ws.Range("O6").Value = 0.6444737471070977
$ws.Range("P6").Value = 0.6444737471070977
$ws.Range("S6").Value = 0.08476821519564612
$ws.Range("T6").Value = 0.08476821519564612

$ws.Range("I7").Value = 0.1315309049843414
$ws.Range("J7").Value = 0.1315309049843414
$ws.Range("O7").Value = 0.09604564432568881
$ws.Range("P7").Value = 0.09604564432568881
$ws.Range("S7").Value = 0.01263297051796202
$ws.Range("T7").Value = 0.01263297051796202

$ws.Range("H8").Value = 461.8238680000001
$ws.Range("I8").Value = 0.8053908704318941
$ws.Range("J8").Value = 0.8053908704318941
$ws.Range("M8").Value = 11.146846
$ws.Range("N8").Value = 33.440538
$ws.Range("O8").Value = 0.2594806085672136
$ws.Range("P8").Value = 0.2594806085672136
$ws.Range("Q8").Value = 1715.95984524011
$ws.Range("R8").Value = 15443.63860716099
$ws.Range("S8").Value = 0.2089833131941457
$ws.Range("T8").Value = 0.2089833131941457

$ws.Range("H9").Value = 461.8238680000001
$ws.Range("I9").Value = 0.8053908704318941
$ws.Range("J9").Value = 0.8053908704318941
$ws.Range("O9").Value = 0.6444737471070977
$ws.Range("P9").Value = 0.6444737471070977
$ws.Range("S9").Value = 0.5190532721530898
$ws.Range("T9").Value = 0.5190532721530898

$ws.Range("H10").Value = 461.8238680000001
$ws.Range("I10").Value = 0.8053908704318941
$ws.Range("J10").Value = 0.8053908704318941
$ws.Range("O10").Value = 0.09604564432568881
$ws.Range("P10").Value = 0.09604564432568881
$ws.Range("S10").Value = 0.07735428508465862
$ws.Range("T10").Value = 0.07735428508465862
